$d = $word.ActiveDocument

$d.Content.Find.Execute("it is a essential to have both.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "it is essential to have both.", 2)
